# Add production information for 3 new WW products (plus a TAS variant row)
# to the label library worksheet. Four new data rows are inserted directly
# above the last ("TEST") row, which shifts from row 159 down to row 163.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

$insertAt = 159

# Insert four blank rows above the existing last row (currently row 159,
# the "TEST" sentinel row). After this, that sentinel row becomes row 163.
for ($i = 0; $i -lt 4; $i++) {
    $ws.Rows.Item($insertAt).Insert()
}

$lastRow = 163

# Copy the (now shifted) sentinel row's formatting into the freshly
# inserted, still-blank rows so the new rows pick up the same cell
# styles / row height used throughout the table.
$ws.Range("A" + $lastRow + ":I" + $lastRow).Copy()
$ws.Range("A159:I162").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# PasteSpecial only carries cell-level formatting, not the row height,
# so match the standard 17.25pt row height used by every other data row.
$ws.Range("A159:I162").RowHeight = 17.25

# NOTE: this interpreter does not reliably bind *named* function
# parameters, so the helper below takes plain positional parameters.
function Set-ProductRow {
    param($Row, $ProductName, $TopLabel, $CartonLabel, $PalletLabel, $WaterMark, $ShelfLife, $LastCol)

    $ws.Range("A" + $Row).Value2 = $ProductName
    $ws.Range("B" + $Row).Value2 = $TopLabel
    $ws.Range("E" + $Row).Value2 = $CartonLabel
    $ws.Range("F" + $Row).Value2 = $PalletLabel
    $ws.Range("G" + $Row).Value2 = $WaterMark
    $ws.Range("H" + $Row).Value2 = $ShelfLife

    # Column I already carries a Text ("@") number format inherited from
    # the copied row; temporarily switch to General so the numeric value
    # is stored as a real number (matching the rest of the sheet), then
    # restore the Text format for display/style consistency.
    $ws.Range("I" + $Row).NumberFormat = "General"
    $ws.Range("I" + $Row).Value2 = $LastCol
    $ws.Range("I" + $Row).NumberFormat = "@"
}

Set-ProductRow 159 "WW Mexican Style Slaw 250gx6" "9339687425990" "19339687425997" "563290" "T4K" "12" 28

Set-ProductRow 160 "WW Mexican Style Slaw 250gx6 (TAS)" "9339687425990" "19339687425997" "563290T" "T4K" "12" 28

Set-ProductRow 161 "WW BBQ Potato Salad 400gx8" "9339687426010" "19339687426017" "563291" "3TK" "12" 28

Set-ProductRow 162 "WW Creamy Rench Pasta Salad 350gx8" "9339687426003" "19339687426000" "563289" "FV6" "12" 28

# Extend the hidden _FilterDatabase named range so it still spans column I
# of the whole table (was I1:I159, now I1:I163).
$filterName = $wb.Names.Item("Sheet1!_FilterDatabase")
$filterName.RefersTo = "=Sheet1!`$I`$1:`$I`$" + $lastRow

# Update the worksheet's active cell / selection to reflect the new last
# row, mirroring where the author ended up after the edit.
$ws.Range("G" + $lastRow).Select()

Write-Output "Inserted 4 product rows (159-162); sentinel row now at 163."
